$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Rg22"
$ws.Range("B5").Value = "Cam"
$ws.Range("C5").Value = "High locality at line speed of 30mpm. Set to 20 dropped it down to passing levels"
$ws.Range("D5").Value = "2025-10-13 12:01:44"
